# Add a new column "CO2/(CO+CO2)" as column H (between "feed Yh2" and "Yco"),
# containing the formula =F/(E+F) for each data row, and shift the existing
# H:N columns (Yco ... wcat (g)) one column to the right (I:O).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at H; this shifts H:Q -> I:R (incl. column
# widths, the used-range dimension and the empty Q2:Q19 spacer cells).
$ws.Columns.Item(8).Insert()

# Header for the new column.
$ws.Range("H1").Value = "CO2/(CO+CO2)"

# Data rows: 2 through 19. Row 2 gets its own formula; rows 3:19 are filled
# as one range assignment (Excel authors this as a shared formula group).
$ws.Range("H2").Formula = "=F2/(E2+F2)"
$ws.Range("H3:H19").Formula = "=F3/(E3+F3)"

# Match the author's final selection (cell H5).
$ws.Range("H5").Select() | Out-Null
